$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 "Categoría" reusing the same header style as A1:E1
$ws.Range("A1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "Categoría"

# Row 2: Salchipapa
$ws.Range("A2").Value = "Salchipapa"
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 2.5
$ws.Range("E2").Value = 1.5
$ws.Range("F2").Value = "Comida"

# Row 3: Coca cola
$ws.Range("A3").Value = "Coca cola"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = "Bebida"

# Row 4: Gorro de baño
$ws.Range("A4").Value = "Gorro de baño"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "Otros"

# Row 5: Chochos con tostado
$ws.Range("A5").Value = "Chochos con tostado"
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 0.25
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 0.75
$ws.Range("F5").Value = "Comida"
